$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 235; this shifts existing rows 235:292 down to 236:293
$ws.Rows.Item(235).Insert()

# Populate the new row 235 with the new data record
$ws.Cells.Item(235, 1).Value = 10
$ws.Cells.Item(235, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(235, 3).Value = "La Araucanía"
$ws.Cells.Item(235, 4).Value = 44855
$ws.Cells.Item(235, 5).Value = 9
$ws.Cells.Item(235, 6).Value = 100112039
$ws.Cells.Item(235, 7).Value = "Ciboulette"
$ws.Cells.Item(235, 8).Value = "Sin especificar"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 50
$ws.Cells.Item(235, 11).Value = 7000
$ws.Cells.Item(235, 12).Value = 8000
$ws.Cells.Item(235, 13).Value = 7400
$ws.Cells.Item(235, 14).Value = "$/docena de atados"
$ws.Cells.Item(235, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(235, 16).Value = 2467
$ws.Cells.Item(235, 17).Value = 3
$ws.Cells.Item(235, 18).Value = "Hortaliza"
